$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the previous merges on row 3 (A3:B3 and B3:C3) first, so that the
# individual header cells on that row can be set independently
$ws.Range("A3:B3").UnMerge()
$ws.Range("B3:C3").UnMerge()

# Update the "Fecha" (date) cell in the title block
$ws.Range("D1").Value = "Fecha  2023-02-04 22:13:24"

# Add the new headers in row 3 (C3:G3), matching style of existing header cells
$ws.Range("C3").Value = "Adicional"
$ws.Range("D3").Value = "Precio"
$ws.Range("E3").Value = "Tipo"
$ws.Range("F3").Value = "Presentacion"
$ws.Range("G3").Value = "Laboratorio"

$ws.Range("A3:G3").Font.Bold = $true
$ws.Range("A3:G3").Font.Size = 12

# Update product data
$ws.Range("A4").Value = "buscapina "
$ws.Range("D5").Value = 5
$ws.Range("F5").Value = "Suero"
$ws.Range("F7").Value = "Suero"
